$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 608.1667
$ws.Cells.Item(6, 9).Value = 430
$ws.Cells.Item(6, 10).Value = 1499
$ws.Cells.Item(6, 11).Value = 1290
$ws.Cells.Item(6, 12).Value = 4497
$ws.Cells.Item(6, 13).Value = -1178
$ws.Cells.Item(6, 14).Value = -4721
$ws.Cells.Item(12, 8).Value = 13902221
$ws.Cells.Item(12, 9).Value = 13902221
$ws.Cells.Item(12, 11).Value = 13902221
$ws.Cells.Item(12, 13).Value = -13902051
$ws.Cells.Item(17, 8).Value = 16667276
$ws.Cells.Item(17, 10).Value = 16667276
$ws.Cells.Item(17, 12).Value = 50001828
$ws.Cells.Item(17, 14).Value = -50002164
$ws.Cells.Item(28, 8).Value = 262.5
$ws.Cells.Item(28, 9).Value = 187.5
$ws.Cells.Item(28, 11).Value = 187.5
$ws.Cells.Item(28, 13).Value = 297.5
$ws.Cells.Item(33, 8).Value = 8936711
$ws.Cells.Item(33, 9).Value = 9657.048000000001
$ws.Cells.Item(33, 11).Value = 9657.048000000001
$ws.Cells.Item(33, 13).Value = -9428.048000000001
$ws.Cells.Item(40, 8).Value = 2523.5
$ws.Cells.Item(40, 10).Value = 2562
$ws.Cells.Item(40, 12).Value = 2562
$ws.Cells.Item(40, 14).Value = -2912
$ws.Cells.Item(53, 8).Value = 71429230
$ws.Cells.Item(53, 9).Value = 192
$ws.Cells.Item(53, 10).Value = 125001016
$ws.Cells.Item(53, 11).Value = 192
$ws.Cells.Item(53, 12).Value = 125001016
$ws.Cells.Item(53, 13).Value = 445
$ws.Cells.Item(53, 14).Value = -125002290
$ws.Cells.Item(62, 8).Value = 4467.28
$ws.Cells.Item(62, 9).Value = 4120.6924
$ws.Cells.Item(62, 10).Value = 4842.75
$ws.Cells.Item(62, 11).Value = 4120.6924
$ws.Cells.Item(62, 12).Value = 4842.75
$ws.Cells.Item(62, 13).Value = -3496.6924
$ws.Cells.Item(62, 14).Value = -6090.75
$ws.Cells.Item(65, 8).Value = 4467.28
$ws.Cells.Item(65, 9).Value = 4120.6924
$ws.Cells.Item(65, 10).Value = 4842.75
$ws.Cells.Item(65, 11).Value = 20603.462
$ws.Cells.Item(65, 12).Value = 24213.75
$ws.Cells.Item(65, 13).Value = -17483.462
$ws.Cells.Item(65, 14).Value = -30453.75
$ws.Cells.Item(70, 8).Value = 2268.6
$ws.Cells.Item(70, 9).Value = 1499
$ws.Cells.Item(70, 10).Value = 2461
$ws.Cells.Item(70, 11).Value = 4497
$ws.Cells.Item(70, 12).Value = 7383
$ws.Cells.Item(70, 13).Value = -4227
$ws.Cells.Item(70, 14).Value = -7923
$ws.Cells.Item(73, 8).Value = 2268.6
$ws.Cells.Item(73, 9).Value = 1499
$ws.Cells.Item(73, 10).Value = 2461
$ws.Cells.Item(73, 11).Value = 4497
$ws.Cells.Item(73, 12).Value = 7383
$ws.Cells.Item(73, 13).Value = -3561
$ws.Cells.Item(73, 14).Value = -9255
$ws.Cells.Item(88, 8).Value = 10044
$ws.Cells.Item(88, 10).Value = 10740.333
$ws.Cells.Item(88, 12).Value = 10740.333
$ws.Cells.Item(88, 14).Value = -11552.333
$ws.Cells.Item(91, 8).Value = 10044
$ws.Cells.Item(91, 10).Value = 10740.333
$ws.Cells.Item(91, 12).Value = 10740.333
$ws.Cells.Item(91, 14).Value = -13548.333
$ws.Cells.Item(92, 8).Value = 2468165
$ws.Cells.Item(92, 9).Value = 1117224.6
$ws.Cells.Item(92, 10).Value = 6250798
$ws.Cells.Item(92, 11).Value = 1117224.6
$ws.Cells.Item(92, 12).Value = 6250798
$ws.Cells.Item(92, 13).Value = -1115976.6
$ws.Cells.Item(92, 14).Value = -6253294
$ws.Cells.Item(97, 8).Value = 3668.9092
$ws.Cells.Item(97, 10).Value = 3668.9092
$ws.Cells.Item(97, 12).Value = 11006.7276
$ws.Cells.Item(97, 14).Value = -11998.7276
$ws.Cells.Item(100, 8).Value = 5013.364
$ws.Cells.Item(100, 9).Value = 1622.5
$ws.Cells.Item(100, 11).Value = 1622.5
$ws.Cells.Item(100, 13).Value = -1081.5
$ws.Cells.Item(101, 8).Value = 1864
$ws.Cells.Item(101, 10).Value = 2395
$ws.Cells.Item(101, 12).Value = 7185
$ws.Cells.Item(101, 14).Value = -10429
$ws.Cells.Item(104, 8).Value = 1200
$ws.Cells.Item(104, 9).Value = 1000
$ws.Cells.Item(104, 10).Value = 1600
$ws.Cells.Item(104, 11).Value = 3000
$ws.Cells.Item(104, 12).Value = 4800
$ws.Cells.Item(104, 13).Value = -1253
$ws.Cells.Item(104, 14).Value = -8294
$ws.Cells.Item(106, 8).Value = 103450744
$ws.Cells.Item(106, 9).Value = 130436990
$ws.Cells.Item(106, 10).Value = 3433.1667
$ws.Cells.Item(106, 11).Value = 130436990
$ws.Cells.Item(106, 12).Value = 3433.1667
$ws.Cells.Item(106, 13).Value = -130436359
$ws.Cells.Item(106, 14).Value = -4695.1667
$ws.Cells.Item(112, 8).Value = 2480.4285
$ws.Cells.Item(112, 9).Value = 1340
$ws.Cells.Item(112, 11).Value = 4020
$ws.Cells.Item(112, 13).Value = -2912
$ws.Cells.Item(113, 8).Value = 3290.4167
$ws.Cells.Item(113, 10).Value = 3166.5557
$ws.Cells.Item(113, 12).Value = 3166.5557
$ws.Cells.Item(113, 14).Value = -9674.555700000001
$ws.Cells.Item(125, 8).Value = 3265
$ws.Cells.Item(125, 9).Value = 2397.5
$ws.Cells.Item(125, 11).Value = 21577.5
$ws.Cells.Item(125, 13).Value = -19117.5
$ws.Cells.Item(132, 8).Value = 1276.6511
$ws.Cells.Item(132, 9).Value = 954.6857
$ws.Cells.Item(132, 11).Value = 2864.0571
$ws.Cells.Item(132, 13).Value = -334.0571
$ws.Cells.Item(137, 8).Value = 2179798.5
$ws.Cells.Item(137, 9).Value = 5562.875
$ws.Cells.Item(137, 10).Value = 4551691.5
$ws.Cells.Item(137, 11).Value = 16688.625
$ws.Cells.Item(137, 12).Value = 13655074.5
$ws.Cells.Item(137, 13).Value = -14138.625
$ws.Cells.Item(137, 14).Value = -13660174.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 820.6316
$ws.Cells.Item(2, 9).Value = 784.5
$ws.Cells.Item(2, 11).Value = 784.5
$ws.Cells.Item(2, 13).Value = -671.5
$ws.Cells.Item(32, 8).Value = 5258.591
$ws.Cells.Item(32, 9).Value = 2289.0356
$ws.Cells.Item(32, 10).Value = 10455.3125
$ws.Cells.Item(32, 11).Value = 2289.0356
$ws.Cells.Item(32, 12).Value = 10455.3125
$ws.Cells.Item(32, 13).Value = -2002.0356
$ws.Cells.Item(32, 14).Value = -11029.3125
$ws.Cells.Item(45, 8).Value = 64684.125
$ws.Cells.Item(45, 9).Value = 64684.125
$ws.Cells.Item(45, 11).Value = 64684.125
$ws.Cells.Item(45, 13).Value = -64307.125
$ws.Cells.Item(61, 8).Value = 2398418
$ws.Cells.Item(61, 9).Value = 77221.14
$ws.Cells.Item(61, 10).Value = 5352668.5
$ws.Cells.Item(61, 11).Value = 77221.14
$ws.Cells.Item(61, 12).Value = 5352668.5
$ws.Cells.Item(61, 13).Value = -77009.14
$ws.Cells.Item(61, 14).Value = -5353092.5
$ws.Cells.Item(97, 8).Value = 7144
$ws.Cells.Item(97, 9).Value = 8658.691999999999
$ws.Cells.Item(97, 10).Value = 2221.25
$ws.Cells.Item(97, 11).Value = 8658.691999999999
$ws.Cells.Item(97, 12).Value = 2221.25
$ws.Cells.Item(97, 13).Value = -8162.691999999999
$ws.Cells.Item(97, 14).Value = -3213.25
$ws.Cells.Item(102, 8).Value = 9909.833000000001
$ws.Cells.Item(102, 9).Value = 11491.8
$ws.Cells.Item(102, 11).Value = 11491.8
$ws.Cells.Item(102, 13).Value = -9869.799999999999
$ws.Cells.Item(107, 8).Value = 55000
$ws.Cells.Item(107, 10).Value = 55000
$ws.Cells.Item(107, 12).Value = 55000
$ws.Cells.Item(107, 14).Value = -62680
$ws.Cells.Item(109, 8).Value = 79750
$ws.Cells.Item(109, 10).Value = 79750
$ws.Cells.Item(109, 12).Value = 79750
$ws.Cells.Item(109, 14).Value = -82524
$ws.Cells.Item(116, 8).Value = 820.6316
$ws.Cells.Item(116, 9).Value = 784.5
$ws.Cells.Item(116, 11).Value = 784.5
$ws.Cells.Item(116, 13).Value = 1509.5
$ws.Cells.Item(122, 8).Value = 2070.0967
$ws.Cells.Item(122, 9).Value = 1945.3462
$ws.Cells.Item(122, 11).Value = 5836.0386
$ws.Cells.Item(122, 13).Value = -3386.0386
$ws.Cells.Item(124, 8).Value = 27944
$ws.Cells.Item(124, 10).Value = 27944
$ws.Cells.Item(124, 12).Value = 27944
$ws.Cells.Item(124, 14).Value = -37764
$ws.Cells.Item(132, 8).Value = 1958.3125
$ws.Cells.Item(132, 9).Value = 1619.1714
$ws.Cells.Item(132, 10).Value = 2871.3845
$ws.Cells.Item(132, 11).Value = 4857.5142
$ws.Cells.Item(132, 12).Value = 8614.1535
$ws.Cells.Item(132, 13).Value = -2327.5142
$ws.Cells.Item(132, 14).Value = -13674.1535
$ws.Cells.Item(135, 8).Value = 67311.14
$ws.Cells.Item(135, 10).Value = 67311.14
$ws.Cells.Item(135, 12).Value = 67311.14
$ws.Cells.Item(135, 14).Value = -77451.14
$ws.Cells.Item(136, 8).Value = 2398418
$ws.Cells.Item(136, 9).Value = 77221.14
$ws.Cells.Item(136, 10).Value = 5352668.5
$ws.Cells.Item(136, 11).Value = 231663.42
$ws.Cells.Item(136, 12).Value = 16058005.5
$ws.Cells.Item(136, 13).Value = -229113.42
$ws.Cells.Item(136, 14).Value = -16063105.5
$ws.Cells.Item(138, 8).Value = 98994.5
$ws.Cells.Item(138, 9).Value = 98489
$ws.Cells.Item(138, 10).Value = 99163
$ws.Cells.Item(138, 11).Value = 98489
$ws.Cells.Item(138, 12).Value = 99163
$ws.Cells.Item(138, 13).Value = -93349
$ws.Cells.Item(138, 14).Value = -109443
$ws.Cells.Item(139, 8).Value = 69377.60000000001
$ws.Cells.Item(139, 10).Value = 69377.60000000001
$ws.Cells.Item(139, 12).Value = 69377.60000000001
$ws.Cells.Item(139, 14).Value = -79657.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 820.6316
$ws.Cells.Item(3, 9).Value = 784.5
$ws.Cells.Item(3, 11).Value = 784.5
$ws.Cells.Item(3, 13).Value = -670.5
$ws.Cells.Item(86, 8).Value = 4512.5864
$ws.Cells.Item(86, 9).Value = 3027.652
$ws.Cells.Item(86, 11).Value = 3027.652
$ws.Cells.Item(86, 13).Value = -1904.652
$ws.Cells.Item(88, 8).Value = 40000
$ws.Cells.Item(88, 10).Value = 40000
$ws.Cells.Item(88, 12).Value = 40000
$ws.Cells.Item(88, 14).Value = -40812
$ws.Cells.Item(89, 8).Value = 4512.5864
$ws.Cells.Item(89, 9).Value = 3027.652
$ws.Cells.Item(89, 11).Value = 15138.26
$ws.Cells.Item(89, 13).Value = -9522.26
$ws.Cells.Item(91, 8).Value = 40000
$ws.Cells.Item(91, 10).Value = 40000
$ws.Cells.Item(91, 12).Value = 40000
$ws.Cells.Item(91, 14).Value = -42808
$ws.Cells.Item(99, 8).Value = 6090.41
$ws.Cells.Item(99, 9).Value = 5610.294
$ws.Cells.Item(99, 11).Value = 5610.294
$ws.Cells.Item(99, 13).Value = -4112.294
$ws.Cells.Item(105, 8).Value = 15827.391
$ws.Cells.Item(105, 9).Value = 20567.035
$ws.Cells.Item(105, 10).Value = 4373.25
$ws.Cells.Item(105, 11).Value = 20567.035
$ws.Cells.Item(105, 12).Value = 4373.25
$ws.Cells.Item(105, 13).Value = -18820.035
$ws.Cells.Item(105, 14).Value = -7867.25
$ws.Cells.Item(107, 8).Value = 18659.77
$ws.Cells.Item(107, 9).Value = 21161.545
$ws.Cells.Item(107, 10).Value = 4900
$ws.Cells.Item(107, 11).Value = 21161.545
$ws.Cells.Item(107, 12).Value = 4900
$ws.Cells.Item(107, 13).Value = -19241.545
$ws.Cells.Item(107, 14).Value = -8740
$ws.Cells.Item(134, 8).Value = 21430842
$ws.Cells.Item(134, 9).Value = 1973.4242
$ws.Cells.Item(134, 10).Value = 100003360
$ws.Cells.Item(134, 11).Value = 5920.2726
$ws.Cells.Item(134, 12).Value = 300010080
$ws.Cells.Item(134, 13).Value = -3385.2726
$ws.Cells.Item(134, 14).Value = -300015150

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 4499.6665
$ws.Cells.Item(6, 10).Value = 10000
$ws.Cells.Item(6, 12).Value = 10000
$ws.Cells.Item(6, 14).Value = -10226
$ws.Cells.Item(7, 8).Value = 295.5625
$ws.Cells.Item(7, 9).Value = 224.44444
$ws.Cells.Item(7, 10).Value = 387
$ws.Cells.Item(7, 11).Value = 224.44444
$ws.Cells.Item(7, 12).Value = 387
$ws.Cells.Item(7, 13).Value = -111.44444
$ws.Cells.Item(7, 14).Value = -613
$ws.Cells.Item(16, 8).Value = 4333649
$ws.Cells.Item(16, 9).Value = 8929750
$ws.Cells.Item(16, 10).Value = 7906.7646
$ws.Cells.Item(16, 11).Value = 8929750
$ws.Cells.Item(16, 12).Value = 7906.7646
$ws.Cells.Item(16, 13).Value = -8929463
$ws.Cells.Item(16, 14).Value = -8480.7646
$ws.Cells.Item(31, 8).Value = 2580.7097
$ws.Cells.Item(31, 9).Value = 4926.75
$ws.Cells.Item(31, 10).Value = 2233.1482
$ws.Cells.Item(31, 11).Value = 4926.75
$ws.Cells.Item(31, 12).Value = 2233.1482
$ws.Cells.Item(31, 13).Value = -4631.75
$ws.Cells.Item(31, 14).Value = -2823.1482
$ws.Cells.Item(34, 8).Value = 2580.7097
$ws.Cells.Item(34, 9).Value = 4926.75
$ws.Cells.Item(34, 10).Value = 2233.1482
$ws.Cells.Item(34, 11).Value = 4926.75
$ws.Cells.Item(34, 12).Value = 2233.1482
$ws.Cells.Item(34, 13).Value = -4724.75
$ws.Cells.Item(34, 14).Value = -2637.1482
$ws.Cells.Item(58, 8).Value = 5774.6924
$ws.Cells.Item(58, 9).Value = 11918.223
$ws.Cells.Item(58, 10).Value = 2522.2354
$ws.Cells.Item(58, 11).Value = 11918.223
$ws.Cells.Item(58, 12).Value = 2522.2354
$ws.Cells.Item(58, 13).Value = -11715.223
$ws.Cells.Item(58, 14).Value = -2928.2354
$ws.Cells.Item(86, 8).Value = 11825.941
$ws.Cells.Item(86, 10).Value = 13360.571
$ws.Cells.Item(86, 12).Value = 13360.571
$ws.Cells.Item(86, 14).Value = -15606.571
$ws.Cells.Item(88, 8).Value = 25848.6
$ws.Cells.Item(88, 10).Value = 25848.6
$ws.Cells.Item(88, 12).Value = 25848.6
$ws.Cells.Item(88, 14).Value = -26660.6
$ws.Cells.Item(89, 8).Value = 11825.941
$ws.Cells.Item(89, 10).Value = 13360.571
$ws.Cells.Item(89, 12).Value = 66802.855
$ws.Cells.Item(89, 14).Value = -78034.855
$ws.Cells.Item(91, 8).Value = 25848.6
$ws.Cells.Item(91, 10).Value = 25848.6
$ws.Cells.Item(91, 12).Value = 25848.6
$ws.Cells.Item(91, 14).Value = -28656.6
$ws.Cells.Item(99, 8).Value = 4657.2856
$ws.Cells.Item(99, 9).Value = 4055.7144
$ws.Cells.Item(99, 11).Value = 4055.7144
$ws.Cells.Item(99, 13).Value = -2557.7144
$ws.Cells.Item(105, 8).Value = 1053
$ws.Cells.Item(105, 9).Value = 1065.3334
$ws.Cells.Item(105, 11).Value = 1065.3334
$ws.Cells.Item(105, 13).Value = 681.6666
$ws.Cells.Item(107, 8).Value = 1124.1471
$ws.Cells.Item(107, 9).Value = 1244.762
$ws.Cells.Item(107, 10).Value = 929.3077
$ws.Cells.Item(107, 11).Value = 1244.762
$ws.Cells.Item(107, 12).Value = 929.3077
$ws.Cells.Item(107, 13).Value = 675.2380000000001
$ws.Cells.Item(107, 14).Value = -4769.3077
$ws.Cells.Item(109, 8).Value = 35583
$ws.Cells.Item(109, 9).Value = 32249.5
$ws.Cells.Item(109, 10).Value = 37249.75
$ws.Cells.Item(109, 11).Value = 32249.5
$ws.Cells.Item(109, 12).Value = 37249.75
$ws.Cells.Item(109, 13).Value = -31209.5
$ws.Cells.Item(109, 14).Value = -39329.75
$ws.Cells.Item(113, 8).Value = 4333649
$ws.Cells.Item(113, 9).Value = 8929750
$ws.Cells.Item(113, 10).Value = 7906.7646
$ws.Cells.Item(113, 11).Value = 8929750
$ws.Cells.Item(113, 12).Value = 7906.7646
$ws.Cells.Item(113, 13).Value = -8927580
$ws.Cells.Item(113, 14).Value = -12246.7646
$ws.Cells.Item(122, 8).Value = 1088.9445
$ws.Cells.Item(122, 9).Value = 1088.9445
$ws.Cells.Item(122, 11).Value = 3266.8335
$ws.Cells.Item(122, 13).Value = -816.8335000000002
$ws.Cells.Item(126, 8).Value = 4657.2856
$ws.Cells.Item(126, 9).Value = 4055.7144
$ws.Cells.Item(126, 11).Value = 12167.1432
$ws.Cells.Item(126, 13).Value = -9697.143199999999
$ws.Cells.Item(132, 8).Value = 19610676
$ws.Cells.Item(132, 9).Value = 2307.9092
$ws.Cells.Item(132, 10).Value = 55559348
$ws.Cells.Item(132, 11).Value = 6923.7276
$ws.Cells.Item(132, 12).Value = 166678044
$ws.Cells.Item(132, 13).Value = -4393.7276
$ws.Cells.Item(132, 14).Value = -166683104
$ws.Cells.Item(136, 8).Value = 5774.6924
$ws.Cells.Item(136, 9).Value = 11918.223
$ws.Cells.Item(136, 10).Value = 2522.2354
$ws.Cells.Item(136, 11).Value = 35754.669
$ws.Cells.Item(136, 12).Value = 7566.706200000001
$ws.Cells.Item(136, 13).Value = -33204.669
$ws.Cells.Item(136, 14).Value = -12666.7062
$ws.Cells.Item(141, 8).Value = 67391.63
$ws.Cells.Item(141, 10).Value = 68789.125
$ws.Cells.Item(141, 12).Value = 68789.125
$ws.Cells.Item(141, 14).Value = -79149.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 468.92307
$ws.Cells.Item(2, 10).Value = 699.875
$ws.Cells.Item(2, 12).Value = 4199.25
$ws.Cells.Item(2, 14).Value = -4425.25
$ws.Cells.Item(12, 8).Value = 42579.918
$ws.Cells.Item(12, 9).Value = 1787
$ws.Cells.Item(12, 10).Value = 46288.363
$ws.Cells.Item(12, 11).Value = 5361
$ws.Cells.Item(12, 12).Value = 138865.089
$ws.Cells.Item(12, 13).Value = -5188
$ws.Cells.Item(12, 14).Value = -139211.089
$ws.Cells.Item(14, 8).Value = 1336.2142
$ws.Cells.Item(14, 9).Value = 1336.2142
$ws.Cells.Item(14, 11).Value = 4008.6426
$ws.Cells.Item(14, 13).Value = -3835.6426
$ws.Cells.Item(38, 8).Value = 375.875
$ws.Cells.Item(38, 9).Value = 482
$ws.Cells.Item(38, 11).Value = 1446
$ws.Cells.Item(38, 13).Value = -1099
$ws.Cells.Item(62, 8).Value = 8022
$ws.Cells.Item(62, 10).Value = 8399.75
$ws.Cells.Item(62, 12).Value = 25199.25
$ws.Cells.Item(62, 14).Value = -26571.25
$ws.Cells.Item(63, 8).Value = 9810.147000000001
$ws.Cells.Item(63, 9).Value = 3260
$ws.Cells.Item(63, 11).Value = 9780
$ws.Cells.Item(63, 13).Value = -9031
$ws.Cells.Item(65, 8).Value = 8022
$ws.Cells.Item(65, 10).Value = 8399.75
$ws.Cells.Item(65, 12).Value = 75597.75
$ws.Cells.Item(65, 14).Value = -82461.75
$ws.Cells.Item(66, 8).Value = 9810.147000000001
$ws.Cells.Item(66, 9).Value = 3260
$ws.Cells.Item(66, 11).Value = 29340
$ws.Cells.Item(66, 13).Value = -25596
$ws.Cells.Item(68, 8).Value = 2048
$ws.Cells.Item(68, 10).Value = 2072
$ws.Cells.Item(68, 12).Value = 6216
$ws.Cells.Item(68, 14).Value = -7838
$ws.Cells.Item(71, 8).Value = 2048
$ws.Cells.Item(71, 10).Value = 2072
$ws.Cells.Item(71, 12).Value = 18648
$ws.Cells.Item(71, 14).Value = -26760
$ws.Cells.Item(99, 8).Value = 49298.477
$ws.Cells.Item(99, 9).Value = 146282.72
$ws.Cells.Item(99, 10).Value = 6867.875
$ws.Cells.Item(99, 11).Value = 438848.16
$ws.Cells.Item(99, 12).Value = 20603.625
$ws.Cells.Item(99, 13).Value = -436602.16
$ws.Cells.Item(99, 14).Value = -25095.625
$ws.Cells.Item(104, 8).Value = 7083.3335
$ws.Cells.Item(104, 10).Value = 7083.3335
$ws.Cells.Item(104, 12).Value = 21250.0005
$ws.Cells.Item(104, 14).Value = -26492.0005
$ws.Cells.Item(118, 8).Value = 26500
$ws.Cells.Item(118, 9).Value = 34333.332
$ws.Cells.Item(118, 10).Value = 3000
$ws.Cells.Item(118, 11).Value = 102999.996
$ws.Cells.Item(118, 12).Value = 9000
$ws.Cells.Item(118, 13).Value = -101756.996
$ws.Cells.Item(118, 14).Value = -11486
$ws.Cells.Item(122, 8).Value = 4762123.5
$ws.Cells.Item(122, 9).Value = 427.4
$ws.Cells.Item(122, 11).Value = 3846.6
$ws.Cells.Item(122, 13).Value = -1396.6
$ws.Cells.Item(130, 8).Value = 13333
$ws.Cells.Item(130, 9).Value = 9999.5
$ws.Cells.Item(130, 10).Value = 20000
$ws.Cells.Item(130, 11).Value = 29998.5
$ws.Cells.Item(130, 12).Value = 60000
$ws.Cells.Item(130, 13).Value = -24978.5
$ws.Cells.Item(130, 14).Value = -70040
$ws.Cells.Item(132, 8).Value = 78118.53999999999
$ws.Cells.Item(132, 9).Value = 1015.6
$ws.Cells.Item(132, 10).Value = 126307.875
$ws.Cells.Item(132, 11).Value = 9140.4
$ws.Cells.Item(132, 12).Value = 1136770.875
$ws.Cells.Item(132, 13).Value = -6610.4
$ws.Cells.Item(132, 14).Value = -1141830.875
$ws.Cells.Item(137, 8).Value = 2265.7
$ws.Cells.Item(137, 9).Value = 951.1429000000001
$ws.Cells.Item(137, 11).Value = 2853.4287
$ws.Cells.Item(137, 13).Value = 2246.5713

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 15555
$ws.Cells.Item(80, 9).Value = 13217.789
$ws.Cells.Item(80, 11).Value = 13217.789
$ws.Cells.Item(80, 13).Value = -12219.789
$ws.Cells.Item(83, 8).Value = 15555
$ws.Cells.Item(83, 9).Value = 13217.789
$ws.Cells.Item(83, 11).Value = 66088.94500000001
$ws.Cells.Item(83, 13).Value = -61096.94500000001
$ws.Cells.Item(97, 8).Value = 745.5625
$ws.Cells.Item(97, 9).Value = 670.1
$ws.Cells.Item(97, 11).Value = 670.1
$ws.Cells.Item(97, 13).Value = -174.1
$ws.Cells.Item(102, 8).Value = 12196934
$ws.Cells.Item(102, 9).Value = 13515327
$ws.Cells.Item(102, 11).Value = 13515327
$ws.Cells.Item(102, 13).Value = -13513705
$ws.Cells.Item(113, 8).Value = 2599.45
$ws.Cells.Item(113, 9).Value = 2207.2144
$ws.Cells.Item(113, 10).Value = 3514.6667
$ws.Cells.Item(113, 11).Value = 2207.2144
$ws.Cells.Item(113, 12).Value = 3514.6667
$ws.Cells.Item(113, 13).Value = -37.21439999999984
$ws.Cells.Item(113, 14).Value = -7854.6667
$ws.Cells.Item(121, 8).Value = 24844
$ws.Cells.Item(121, 10).Value = 24844
$ws.Cells.Item(121, 12).Value = 24844
$ws.Cells.Item(121, 14).Value = -28338
$ws.Cells.Item(122, 8).Value = 2725.5715
$ws.Cells.Item(122, 9).Value = 2902.0715
$ws.Cells.Item(122, 11).Value = 8706.2145
$ws.Cells.Item(122, 13).Value = -6256.2145
$ws.Cells.Item(126, 8).Value = 2246.7778
$ws.Cells.Item(126, 9).Value = 1643.8
$ws.Cells.Item(126, 10).Value = 3000.5
$ws.Cells.Item(126, 11).Value = 4931.4
$ws.Cells.Item(126, 12).Value = 9001.5
$ws.Cells.Item(126, 13).Value = -2461.4
$ws.Cells.Item(126, 14).Value = -13941.5
$ws.Cells.Item(132, 8).Value = 7073094.5
$ws.Cells.Item(132, 9).Value = 6371.9
$ws.Cells.Item(132, 11).Value = 19115.7
$ws.Cells.Item(132, 13).Value = -16585.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6906.731
$ws.Cells.Item(7, 9).Value = 3662.2632
$ws.Cells.Item(7, 10).Value = 15713.143
$ws.Cells.Item(7, 11).Value = 3662.2632
$ws.Cells.Item(7, 12).Value = 15713.143
$ws.Cells.Item(7, 13).Value = -3550.2632
$ws.Cells.Item(7, 14).Value = -15937.143
$ws.Cells.Item(16, 8).Value = 3451.7144
$ws.Cells.Item(16, 9).Value = 3193.8333
$ws.Cells.Item(16, 10).Value = 4999
$ws.Cells.Item(16, 11).Value = 3193.8333
$ws.Cells.Item(16, 12).Value = 4999
$ws.Cells.Item(16, 13).Value = -3023.8333
$ws.Cells.Item(16, 14).Value = -5339
$ws.Cells.Item(40, 8).Value = 3710.4285
$ws.Cells.Item(40, 9).Value = 1999.3334
$ws.Cells.Item(40, 10).Value = 4993.75
$ws.Cells.Item(40, 11).Value = 1999.3334
$ws.Cells.Item(40, 12).Value = 4993.75
$ws.Cells.Item(40, 13).Value = -1863.3334
$ws.Cells.Item(40, 14).Value = -5265.75
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 14).ClearContents()
$ws.Cells.Item(60, 8).Value = 83000
$ws.Cells.Item(60, 10).Value = 83000
$ws.Cells.Item(60, 12).Value = 83000
$ws.Cells.Item(60, 14).Value = -84018
$ws.Cells.Item(68, 8).Value = 2720.2222
$ws.Cells.Item(68, 9).Value = 2677
$ws.Cells.Item(68, 11).Value = 2677
$ws.Cells.Item(68, 13).Value = -1928
$ws.Cells.Item(71, 8).Value = 2720.2222
$ws.Cells.Item(71, 9).Value = 2677
$ws.Cells.Item(71, 11).Value = 13385
$ws.Cells.Item(71, 13).Value = -9641
$ws.Cells.Item(82, 8).Value = 2420.6
$ws.Cells.Item(82, 9).Value = 2553.4546
$ws.Cells.Item(82, 10).Value = 2055.25
$ws.Cells.Item(82, 11).Value = 2553.4546
$ws.Cells.Item(82, 12).Value = 2055.25
$ws.Cells.Item(82, 13).Value = -2192.4546
$ws.Cells.Item(82, 14).Value = -2777.25
$ws.Cells.Item(85, 8).Value = 2420.6
$ws.Cells.Item(85, 9).Value = 2553.4546
$ws.Cells.Item(85, 10).Value = 2055.25
$ws.Cells.Item(85, 11).Value = 2553.4546
$ws.Cells.Item(85, 12).Value = 2055.25
$ws.Cells.Item(85, 13).Value = -1305.4546
$ws.Cells.Item(85, 14).Value = -4551.25
$ws.Cells.Item(99, 8).Value = 15416.667
$ws.Cells.Item(99, 9).Value = 15416.667
$ws.Cells.Item(99, 11).Value = 15416.667
$ws.Cells.Item(99, 13).Value = -12421.667
$ws.Cells.Item(100, 8).Value = 4543.6924
$ws.Cells.Item(100, 9).Value = 3556.2856
$ws.Cells.Item(100, 11).Value = 3556.2856
$ws.Cells.Item(100, 13).Value = -3015.2856
$ws.Cells.Item(103, 8).Value = 39982
$ws.Cells.Item(103, 10).Value = 39982
$ws.Cells.Item(103, 12).Value = 39982
$ws.Cells.Item(103, 14).Value = -42326
$ws.Cells.Item(122, 8).Value = 3295.2593
$ws.Cells.Item(122, 9).Value = 2999.913
$ws.Cells.Item(122, 10).Value = 4993.5
$ws.Cells.Item(122, 11).Value = 8999.739
$ws.Cells.Item(122, 12).Value = 14980.5
$ws.Cells.Item(122, 13).Value = -6549.739
$ws.Cells.Item(122, 14).Value = -19880.5
$ws.Cells.Item(126, 8).Value = 6906.731
$ws.Cells.Item(126, 9).Value = 3662.2632
$ws.Cells.Item(126, 10).Value = 15713.143
$ws.Cells.Item(126, 11).Value = 10986.7896
$ws.Cells.Item(126, 12).Value = 47139.429
$ws.Cells.Item(126, 13).Value = -8516.7896
$ws.Cells.Item(126, 14).Value = -52079.429
$ws.Cells.Item(132, 8).Value = 3571.5
$ws.Cells.Item(132, 9).Value = 3335.9048
$ws.Cells.Item(132, 10).Value = 3862.5293
$ws.Cells.Item(132, 11).Value = 10007.7144
$ws.Cells.Item(132, 12).Value = 11587.5879
$ws.Cells.Item(132, 13).Value = -7477.714399999999
$ws.Cells.Item(132, 14).Value = -16647.5879
$ws.Cells.Item(135, 8).Value = 179990
$ws.Cells.Item(135, 10).Value = 179990
$ws.Cells.Item(135, 12).Value = 179990
$ws.Cells.Item(135, 14).Value = -190130
$ws.Cells.Item(136, 8).Value = 2290.0476
$ws.Cells.Item(136, 9).Value = 1474.5
$ws.Cells.Item(136, 10).Value = 3377.4443
$ws.Cells.Item(136, 11).Value = 4423.5
$ws.Cells.Item(136, 12).Value = 10132.3329
$ws.Cells.Item(136, 13).Value = -1873.5
$ws.Cells.Item(136, 14).Value = -15232.3329
$ws.Cells.Item(140, 8).Value = 102666.664
$ws.Cells.Item(140, 10).Value = 102666.664
$ws.Cells.Item(140, 12).Value = 102666.664
$ws.Cells.Item(140, 14).Value = -113026.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(22, 8).Value = 2000
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(62, 8).Value = 3477
$ws.Cells.Item(62, 9).Value = 3477
$ws.Cells.Item(62, 11).Value = 3477
$ws.Cells.Item(62, 13).Value = -2853
$ws.Cells.Item(65, 8).Value = 3477
$ws.Cells.Item(65, 9).Value = 3477
$ws.Cells.Item(65, 11).Value = 17385
$ws.Cells.Item(65, 13).Value = -14265
$ws.Cells.Item(113, 8).Value = 1295
$ws.Cells.Item(113, 10).Value = 687.75
$ws.Cells.Item(113, 12).Value = 2063.25
$ws.Cells.Item(113, 14).Value = -6403.25
$ws.Cells.Item(122, 8).Value = 2189.9614
$ws.Cells.Item(122, 9).Value = 1877.56
$ws.Cells.Item(122, 11).Value = 5632.68
$ws.Cells.Item(122, 13).Value = -3182.68
$ws.Cells.Item(125, 8).Value = 84992.8
$ws.Cells.Item(125, 10).Value = 84992.8
$ws.Cells.Item(125, 12).Value = 84992.8
$ws.Cells.Item(125, 14).Value = -94832.8
$ws.Cells.Item(132, 8).Value = 1601.7142
$ws.Cells.Item(132, 9).Value = 1103.6923
$ws.Cells.Item(132, 10).Value = 2411
$ws.Cells.Item(132, 11).Value = 3311.0769
$ws.Cells.Item(132, 12).Value = 7233
$ws.Cells.Item(132, 13).Value = -781.0769
$ws.Cells.Item(132, 14).Value = -12293
$ws.Cells.Item(136, 8).Value = 5104.968
$ws.Cells.Item(136, 9).Value = 6553.4
$ws.Cells.Item(136, 11).Value = 19660.2
$ws.Cells.Item(136, 13).Value = -17110.2
